# Update template's ontology terms (2EXT04_DNA)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "2EXT04_DNA" (main annotation table)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2EXT04_DNA")

# Row 2: RNA (Transcriptomics)
$ws.Range("B2").Value = "RNA (Transcriptomics)"
$ws.Range("C2").Value = "user-specific"
$ws.Range("D2").Value = "user-specific"
$ws.Range("F2").Value = "milligram"
$ws.Range("G2").Value = "UO"
$ws.Range("H2").Value = "http://purl.obolibrary.org/obo/UO_0000022"
$ws.Range("I2").Value = "QIAGEN RNEasy"
$ws.Range("J2").Value = "user-specific"
$ws.Range("K2").Value = "user-specific"
$ws.Range("L2").Value = "QIAGEN RNEasy Buffer 2"
$ws.Range("M2").Value = "user-specific"
$ws.Range("N2").Value = "user-specific"
$ws.Range("O2").Value = 200
$ws.Range("P2").Value = "microliter"
$ws.Range("Q2").Value = "UO"
$ws.Range("R2").Value = "http://purl.obolibrary.org/obo/UO_0000101"

# Row 3: DNA (Genomics)
$ws.Range("B3").Value = "DNA (Genomics)"
$ws.Range("C3").Value = "user-specific"
$ws.Range("D3").Value = "user-specific"
$ws.Range("F3").Value = "milligram"
$ws.Range("G3").Value = "UO"
$ws.Range("H3").Value = "http://purl.obolibrary.org/obo/UO_0000022"
$ws.Range("I3").Value = "PCI method"
$ws.Range("J3").Value = "user-specific"
$ws.Range("K3").Value = "user-specific"
$ws.Range("L3").Value = "phenol:chloroform:isopropanol"
$ws.Range("M3").Value = "user-specific"
$ws.Range("N3").Value = "user-specific"
$ws.Range("P3").Value = "microliter"
$ws.Range("Q3").Value = "UO"
$ws.Range("R3").Value = "http://purl.obolibrary.org/obo/UO_0000101"

# Row 4: Metabolites
$ws.Range("B4").Value = "Metabolites"
$ws.Range("C4").Value = "user-specific"
$ws.Range("D4").Value = "user-specific"
$ws.Range("F4").Value = "milligram"
$ws.Range("G4").Value = "UO"
$ws.Range("H4").Value = "http://purl.obolibrary.org/obo/UO_0000022"
$ws.Range("P4").Value = "microliter"
$ws.Range("Q4").Value = "UO"
$ws.Range("R4").Value = "http://purl.obolibrary.org/obo/UO_0000101"

# Row 5: protein (was "Protein" with no ontology terms, now lower-case "protein"
# annotated against NFDI4PSO)
$ws.Range("B5").Value = "protein"
$ws.Range("C5").Value = "NFDI4PSO"
$ws.Range("D5").Value = "http://purl.obolibrary.org/obo/NFDI4PSO_1000093"
$ws.Range("F5").Value = "milligram"
$ws.Range("G5").Value = "UO"
$ws.Range("H5").Value = "http://purl.obolibrary.org/obo/UO_0000022"
$ws.Range("P5").Value = "microliter"
$ws.Range("Q5").Value = "UO"
$ws.Range("R5").Value = "http://purl.obolibrary.org/obo/UO_0000101"

# Column widths: previously-empty (width 0) hidden helper columns now carry the
# standard "Term Source REF" / "Term Accession Number" widths, and the hidden
# "Unit" helper column narrows to match its sibling unit columns.
$ws.Columns.Item(3).ColumnWidth = 35.75   # Term Source REF (NFDI4PSO:0000012)
$ws.Columns.Item(3).Hidden = $true
$ws.Columns.Item(4).ColumnWidth = 42.75   # Term Accession Number (NFDI4PSO:0000012)
$ws.Columns.Item(4).Hidden = $true
$ws.Columns.Item(6).ColumnWidth = 6.25    # Unit
$ws.Columns.Item(6).Hidden = $true
$ws.Columns.Item(7).ColumnWidth = 35.75   # Term Source REF (NFDI4PSO:0000013)
$ws.Columns.Item(7).Hidden = $true
$ws.Columns.Item(8).ColumnWidth = 42.75   # Term Accession Number (NFDI4PSO:0000013)
$ws.Columns.Item(8).Hidden = $true

# ---------------------------------------------------------------------------
# Sheet "SwateTemplateMetadata" -> bump template version
# ---------------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("SwateTemplateMetadata")
# Leading apostrophe keeps the cell's existing "quote prefix" (text-forced)
# formatting intact, same as the original "1.1.4" entry.
$wsMeta.Range("B3").Value = "'1.1.5"
